$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bernardo's review count for the "Inacio" column goes from 1 to 2
$ws.Range("L4").Value = 2

# 2. "Inacio's code smells" moves from IN PROGRESS (D15) to REVIEWING (F15),
#    picking up its reviewers (Martin, Francisco) in the card title.
$ws.Range("D15").Value = ""
$dst = $ws.Range("F15")
$dst.Value = "Inacio's code smells (Martin, Francisco)"

# Match the formatting already used by cards in the REVIEWING column
# (centered, vertically centered, wrapped text) instead of the blank-cell look.
$src = $ws.Range("F17")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. "Inacio's design patterns" keeps Carlos as reviewer, drops Bernardo
#    (the docx -> txt review file rename dropped Bernardo's own review task).
$ws.Range("F16").Value = "Inacio's design patterns (Carlos)"
